# "se añadió el HV" - add the HV (hypervolume) column to the hyperparameter
# results table, relabel the old "t" header to "t(s)", and fill in the first
# few rows of timing / solution-count data that were already measured.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 12): K was "t", now "t(s)"; L stays "cant_solu";
#      M is a brand-new "hv" column -------------------------------------
$ws.Range("K12").Value = "t(s)"
$ws.Range("M12").Value = "hv"

# Give the new M12 header the same look (bold, centered, filled, bordered)
# as the rest of the header row by copying the format from its neighbour.
$ws.Range("L12").Copy()
$ws.Range("M12").PasteSpecial(-4122)   # xlPasteFormats

# ---- First block of results (rows 13-16) already has timing / solution
#      counts recorded -----------------------------------------------------
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 23

$ws.Range("K14").Value = 1.1
$ws.Range("L14").Value = 30

$ws.Range("K15").Value = 1.1
$ws.Range("L15").Value = 29

$ws.Range("K16").Value = 2.2
$ws.Range("L16").Value = 30

# ---- Bring the new K:M (t(s) / cant_solu / hv) columns into the same
#      visual style as the rest of each data block (F:H) -------------------
$ws.Range("F13").Copy()
$ws.Range("K13:M20").PasteSpecial(-4122)
$ws.Range("K22:M29").PasteSpecial(-4122)
$ws.Range("K31:M39").PasteSpecial(-4122)
$ws.Range("M40").PasteSpecial(-4122)

# Bottom rows of each block (21 and 30) carry the thicker block-ending
# border, so copy that variant separately.
$ws.Range("F21").Copy()
$ws.Range("K21:M21").PasteSpecial(-4122)
$ws.Range("K30:M30").PasteSpecial(-4122)

# ---- Restore the view: scrolled down one row, selection sitting on the
#      newly added M21 cell -------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 8
$win.ScrollColumn = 1
$ws.Range("M21").Select() | Out-Null
